# Rename the four worksheet tabs from underscore-separated names to
# space-separated names:
#   IMPACT_CONFIG      -> IMPACT CONFIG
#   STRATEGY_SELECTOR  -> STRATEGY SELECTOR
#   UPLOAD_READY_ESG   -> UPLOAD READY ESG
#   CROSS_REFERENCE    -> CROSS REFERENCE
#
# Renaming a sheet in Excel normally causes every formula that refers to it
# to be rewritten automatically (and re-quoted, since the new name contains
# a space). The source workbook, however, has a number of formulas that
# refer to these sheets using their old (now stale) unquoted names -- those
# formula strings must be left exactly as-is. Only the sheet tab names
# themselves and the chart series source references should change.
#
# Strategy:
#   1. Remember the old -> new name mapping.
#   2. Rename each sheet (Excel will auto-rewrite dependent formulas).
#   3. Walk every formula cell on every sheet and undo any auto-rewrite by
#      swapping the quoted new name back to the original unquoted old name.
#   4. Fix up the chart's series formulas, which reference
#      'STRATEGY_SELECTOR' but are NOT auto-updated by the rename, so they
#      must be repointed at the new quoted sheet name explicitly.

$wb = $excel.ActiveWorkbook

$renames = @{
    "IMPACT_CONFIG"     = "IMPACT CONFIG";
    "STRATEGY_SELECTOR" = "STRATEGY SELECTOR";
    "UPLOAD_READY_ESG"  = "UPLOAD READY ESG";
    "CROSS_REFERENCE"   = "CROSS REFERENCE"
}

# --- 1 & 2: rename the sheet tabs -----------------------------------------
foreach ($oldName in $renames.Keys) {
    $newName = $renames[$oldName]
    $wb.Worksheets.Item($oldName).Name = $newName
}

# --- 3: restore formula text that Excel auto-rewrote ----------------------
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ($cell.HasFormula) {
            $formula = $cell.Formula
            $original = $formula
            foreach ($oldName in $renames.Keys) {
                $newName = $renames[$oldName]
                $quotedNew = "'" + $newName + "'!"
                $formula = $formula.Replace($quotedNew, $oldName + "!")
            }
            if ($formula -ne $original) {
                $cell.Formula = $formula
            }
        }
    }
}

# --- 4: repoint the chart series at the renamed STRATEGY SELECTOR sheet ---
foreach ($ws in $wb.Worksheets) {
    $chartObjects = $ws.ChartObjects()
    for ($i = 1; $i -le $chartObjects.Count; $i++) {
        $chart = $chartObjects.Item($i).Chart
        $series = $chart.SeriesCollection()
        for ($j = 1; $j -le $series.Count; $j++) {
            $ser = $series.Item($j)
            $serFormula = $ser.Formula
            $serFormula = $serFormula.Replace("STRATEGY_SELECTOR!", "'STRATEGY SELECTOR'!")
            $ser.Formula = $serFormula
        }
    }
}
